$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H2: "testen" -> "test 1 fr" (already styled s="1", just update the text)
$ws.Range("H2").Value = "test 1 fr"

# Copy the formatting of D2 (style s="1") onto the newly translated cells
# so they pick up the existing style instead of synthesizing a new one.
$ws.Range("D2").Copy()

# D3: brand-new cell -> needs style s="1" and text "test 1 fr"
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "test 1 fr"

# H3: "Beispiel" -> "test 1 fr", also gains style s="1"
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").Value = "test 1 fr"

# D4: brand-new cell -> needs style s="1" and text "test 1 fr"
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = "test 1 fr"

# H4: "angeln" -> "test 1 fr", also gains style s="1"
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H4").Value = "test 1 fr"

$excel.CutCopyMode = 0
